$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header row values (row 1)
$ws.Range("A1").Value = "data_path_0"
$ws.Range("B1").Value = "data_path_1"
$ws.Range("C1").Value = "data_path_2"
$ws.Range("D1").Value = "data_path_3"
$ws.Range("E1").Value = "index_0"
$ws.Range("F1").Value = "index_1"
$ws.Range("G1").Value = "index_2"
$ws.Range("H1").Value = "index_3"

# Apply the existing header style (from A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Set data row values (row 2)
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/actionSequence_3-test-data"
$ws.Range("B2").Value = "Data Files/AI-Generated/Common/actionSequence_3-test-data"
$ws.Range("C2").Value = "Data Files/AI-Generated/Common/actionSequence_3-test-data"
$ws.Range("D2").Value = "Data Files/AI-Generated/Common/actionSequence_5-test-data"

# E2:H2 need to hold the text "1" (not the number 1). Stage it in a helper
# cell formatted as Text so Excel doesn't auto-coerce it to a number, copy
# just the value (not the format) into place, then remove the helper cell.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "1"
$ws.Range("Z1").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("F2").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("G2").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# Set column widths. The stored OOXML <col width> is ColumnWidth plus a
# fixed ~0.8333 padding offset, so back that off here to land on the exact
# target widths of 59 (A:D) and 9 (E:H) characters.
$ws.Columns.Item(1).ColumnWidth = 58.166666666666664
$ws.Columns.Item(2).ColumnWidth = 58.166666666666664
$ws.Columns.Item(3).ColumnWidth = 58.166666666666664
$ws.Columns.Item(4).ColumnWidth = 58.166666666666664
$ws.Columns.Item(5).ColumnWidth = 8.166666666666666
$ws.Columns.Item(6).ColumnWidth = 8.166666666666666
$ws.Columns.Item(7).ColumnWidth = 8.166666666666666
$ws.Columns.Item(8).ColumnWidth = 8.166666666666666
